$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report volume number + date range), preserving the
# existing rich-text run styling by just overwriting the full cell text.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/18/2023  Through  12/24/2023"

# ---------------------------------------------------------------------------
# Helper: convert a numeric cell into the "0" / "***.*" text placeholder
# used throughout this sheet for n/a style entries, while reusing the exact
# style (s="14") already applied to the existing placeholder cells. A plain
# Value assignment re-infers a number (and a bare apostrophe-prefix trick
# bakes in an unwanted quote-prefix style), so instead we copy an existing
# "0"/"***.*" placeholder cell's format+value onto the target cell.
# ---------------------------------------------------------------------------
function Set-NaPlaceholder($targetAddr, $sourceAddr) {
    $ws.Range($sourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($sourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4104) | Out-Null   # xlPasteAll (value)
}

# Source placeholder cells already present (and unchanged) in the sheet:
#   C14 -> text "0"    (shared string 20)
#   E14 -> text "***.*" (shared string 21)

# Row 15 (Rape)
Set-NaPlaceholder "D15" "C14"
Set-NaPlaceholder "E15" "E14"
$ws.Range("L15").Value = -40

# Row 16 (Fel. Assault)
Set-NaPlaceholder "D16" "C14"
Set-NaPlaceholder "E16" "E14"
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 72
$ws.Range("K16").Value = -18.181818181818
$ws.Range("L16").Value = 20
$ws.Range("M16").Value = -39.495798319327
$ws.Range("N16").Value = -88.923076923076

# Row 17 (Burglary)
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 71.428571428571
$ws.Range("I17").Value = 103
$ws.Range("J17").Value = 84
$ws.Range("K17").Value = 22.619047619047
$ws.Range("L17").Value = 56.060606060606
$ws.Range("M17").Value = 77.586206896551
$ws.Range("N17").Value = -16.935483870967

# Row 18 (Gr. Larceny)
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -35.714285714285
$ws.Range("I18").Value = 113
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = -5.833333333333
$ws.Range("L18").Value = 22.826086956521
$ws.Range("M18").Value = -15.671641791044
$ws.Range("N18").Value = -91.757840991976

# Row 19 (G.L.A.)
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -57.142857142857
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 48
$ws.Range("I19").Value = 463
$ws.Range("J19").Value = 490
$ws.Range("K19").Value = -5.510204081632
$ws.Range("L19").Value = 24.797843665768
$ws.Range("M19").Value = 20.88772845953
$ws.Range("N19").Value = -54.021847070506

# Row 20 (TOTAL major crime categories before Transit/Housing split)
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -17.647058823529
$ws.Range("I20").Value = 167
$ws.Range("J20").Value = 106
$ws.Range("K20").Value = 57.547169811320
$ws.Range("L20").Value = 142.028985507246
$ws.Range("M20").Value = 68.686868686868
$ws.Range("N20").Value = -94.983478522078

# Row 21 (TOTAL)
$ws.Range("C21").Value = 12
$ws.Range("E21").Value = -29.411764705882
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 69
$ws.Range("H21").Value = 11.594202898550
$ws.Range("I21").Value = 924
$ws.Range("J21").Value = 905
$ws.Range("K21").Value = 2.099447513812
$ws.Range("L21").Value = 38.323353293413
$ws.Range("M21").Value = 15.789473684210
$ws.Range("N21").Value = -85.775862068965

# Row 22 (Transit)
Set-NaPlaceholder "C22" "C14"
Set-NaPlaceholder "D22" "C14"
Set-NaPlaceholder "E22" "E14"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("L22").Value = 47.368421052631

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 43
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = 26.470588235294
$ws.Range("F24").Value = 133
$ws.Range("G24").Value = 167
$ws.Range("H24").Value = -20.359281437125
$ws.Range("I24").Value = 1534
$ws.Range("J24").Value = 1767
$ws.Range("K24").Value = -13.186191284663
$ws.Range("L24").Value = 10.998552821997
$ws.Range("M24").Value = 56.530612244898

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 11
$ws.Range("E25").Value = 266.666666666667
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 29.411764705882
$ws.Range("I25").Value = 245
$ws.Range("J25").Value = 187
$ws.Range("K25").Value = 31.016042780748
$ws.Range("L25").Value = 36.111111111111
$ws.Range("M25").Value = 13.953488372093

# Row 26 (UCR Rape*)
Set-NaPlaceholder "D26" "C14"
Set-NaPlaceholder "E26" "E14"
$ws.Range("L26").Value = -31.578947368421

# Row 27 (Other Sex Crimes)
Set-NaPlaceholder "C27" "C14"

# Row 30 (Hate Crimes)
Set-NaPlaceholder "C30" "C14"
